$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column L header "break_on_off" (new shared string)
$ws.Range("L1").Value = "break_on_off"

# Fill L2:L73 with break_on_off values (1 marks a break trial, else 0)
$ws.Range("L2").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("L19").Value = 1
$ws.Range("L20").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("L37").Value = 1
$ws.Range("L38").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("L54").Value = 1
$ws.Range("L55").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("L73").Value = 0

# Update selection to match the saved view state (row 18 selected)
[void]$ws.Rows.Item(18).Select()

Write-Host "done"
